$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H46").Value = 0
$ws.Range("J46").Value = 0
$ws.Range("L46").Value = 0
$ws.Range("N46").ClearContents()
$ws.Range("H60").Value = 0
$ws.Range("J60").Value = 0
$ws.Range("L60").Value = 0
$ws.Range("N60").ClearContents()
$ws.Range("H106").Value = 1917.1333
$ws.Range("I106").Value = 1839.7858
$ws.Range("J106").Value = 3000
$ws.Range("K106").Value = 1839.7858
$ws.Range("L106").Value = 3000
$ws.Range("M106").Value = -1208.7858
$ws.Range("N106").Value = -4262
$ws.Range("H127").Value = 8961.916999999999
$ws.Range("I127").Value = 6186.25
$ws.Range("K127").Value = 18558.75
$ws.Range("M127").Value = -13598.75
$ws.Range("H135").Value = 9834.286
$ws.Range("I135").Value = 6026.5713
$ws.Range("K135").Value = 54239.14169999999
$ws.Range("M135").Value = -51704.14169999999
$ws.Range("H138").Value = 3998.9866
$ws.Range("J138").Value = 4105.754
$ws.Range("L138").Value = 12317.262
$ws.Range("N138").Value = -22597.262
$ws.Range("H139").Value = 191043.89
$ws.Range("J139").Value = 202424.38
$ws.Range("L139").Value = 202424.38
$ws.Range("N139").Value = -212704.38

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 14902.47
$ws.Range("I32").Value = 7627.968
$ws.Range("J32").Value = 26771.395
$ws.Range("K32").Value = 7627.968
$ws.Range("L32").Value = 26771.395
$ws.Range("M32").Value = -7340.968
$ws.Range("N32").Value = -27345.395
$ws.Range("H61").Value = 4179.7954
$ws.Range("I61").Value = 3894.6128
$ws.Range("K61").Value = 3894.6128
$ws.Range("M61").Value = -3682.6128
$ws.Range("H88").Value = 4297.143
$ws.Range("I88").Value = 2466.3333
$ws.Range("K88").Value = 2466.3333
$ws.Range("M88").Value = -2060.3333
$ws.Range("H91").Value = 4297.143
$ws.Range("I91").Value = 2466.3333
$ws.Range("K91").Value = 2466.3333
$ws.Range("M91").Value = -1062.3333
$ws.Range("H132").Value = 695777.5
$ws.Range("I132").Value = 1013028
$ws.Range("J132").Value = 46855.953
$ws.Range("K132").Value = 3039084
$ws.Range("L132").Value = 140567.859
$ws.Range("M132").Value = -3036554
$ws.Range("N132").Value = -145627.859
$ws.Range("H136").Value = 4179.7954
$ws.Range("I136").Value = 3894.6128
$ws.Range("K136").Value = 11683.8384
$ws.Range("M136").Value = -9133.838400000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H80").Value = 15163649
$ws.Range("I80").Value = 1483
$ws.Range("K80").Value = 1483
$ws.Range("M80").Value = -485
$ws.Range("H83").Value = 15163649
$ws.Range("I83").Value = 1483
$ws.Range("K83").Value = 7415
$ws.Range("M83").Value = -2423
$ws.Range("H86").Value = 8763.9
$ws.Range("I86").Value = 1700.2
$ws.Range("J86").Value = 15827.6
$ws.Range("K86").Value = 1700.2
$ws.Range("L86").Value = 15827.6
$ws.Range("M86").Value = -577.2
$ws.Range("N86").Value = -18073.6
$ws.Range("H89").Value = 8763.9
$ws.Range("I89").Value = 1700.2
$ws.Range("J89").Value = 15827.6
$ws.Range("K89").Value = 8501
$ws.Range("L89").Value = 79138
$ws.Range("M89").Value = -2885
$ws.Range("N89").Value = -90370
$ws.Range("H134").Value = 986478.5600000001
$ws.Range("I134").Value = 1284065.4
$ws.Range("K134").Value = 3852196.2
$ws.Range("M134").Value = -3849661.2

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6278.2764
$ws.Range("J31").Value = 8343.468999999999
$ws.Range("L31").Value = 8343.468999999999
$ws.Range("N31").Value = -8933.468999999999
$ws.Range("H34").Value = 6278.2764
$ws.Range("J34").Value = 8343.468999999999
$ws.Range("L34").Value = 8343.468999999999
$ws.Range("N34").Value = -8747.468999999999
$ws.Range("H51").Value = 41662.332
$ws.Range("J51").Value = 49993.5
$ws.Range("L51").Value = 49993.5
$ws.Range("N51").Value = -51465.5
$ws.Range("H58").Value = 6525.4375
$ws.Range("I58").Value = 4799.7
$ws.Range("J58").Value = 9401.666999999999
$ws.Range("K58").Value = 4799.7
$ws.Range("L58").Value = 9401.666999999999
$ws.Range("M58").Value = -4596.7
$ws.Range("N58").Value = -9807.666999999999
$ws.Range("H59").Value = 69964.664
$ws.Range("J59").Value = 74947
$ws.Range("L59").Value = 74947
$ws.Range("N59").Value = -77237
$ws.Range("H61").Value = 41662.332
$ws.Range("J61").Value = 49993.5
$ws.Range("L61").Value = 49993.5
$ws.Range("N61").Value = -50689.5
$ws.Range("H68").Value = 0
$ws.Range("J68").Value = 0
$ws.Range("L68").Value = 0
$ws.Range("N68").ClearContents()
$ws.Range("H71").Value = 0
$ws.Range("J71").Value = 0
$ws.Range("L71").Value = 0
$ws.Range("N71").ClearContents()
$ws.Range("H132").Value = 7554.6
$ws.Range("I132").Value = 7113.478
$ws.Range("K132").Value = 21340.434
$ws.Range("M132").Value = -18810.434
$ws.Range("H136").Value = 6525.4375
$ws.Range("I136").Value = 4799.7
$ws.Range("J136").Value = 9401.666999999999
$ws.Range("K136").Value = 14399.1
$ws.Range("L136").Value = 28205.001
$ws.Range("M136").Value = -11849.1
$ws.Range("N136").Value = -33305.001
$ws.Range("H141").Value = 192844.81
$ws.Range("J141").Value = 195839.81
$ws.Range("L141").Value = 195839.81
$ws.Range("N141").Value = -206199.81

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1476.3055
$ws.Range("J5").Value = 2136.0557
$ws.Range("L5").Value = 6408.1671
$ws.Range("N5").Value = -6632.1671
$ws.Range("H34").Value = 3613.3333
$ws.Range("J34").Value = 5217.8887
$ws.Range("L34").Value = 15653.6661
$ws.Range("N34").Value = -15821.6661
$ws.Range("H55").Value = 2504124.8
$ws.Range("J55").Value = 3337499.8
$ws.Range("L55").Value = 10012499.4
$ws.Range("N55").Value = -10012853.4
$ws.Range("H94").Value = 11249.8
$ws.Range("I94").Value = 6833.3335
$ws.Range("J94").Value = 13142.571
$ws.Range("K94").Value = 20500.0005
$ws.Range("L94").Value = 39427.713
$ws.Range("M94").Value = -19824.0005
$ws.Range("N94").Value = -40779.713
$ws.Range("H135").Value = 1476.3055
$ws.Range("J135").Value = 2136.0557
$ws.Range("L135").Value = 19224.5013
$ws.Range("N135").Value = -24294.5013

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H44").Value = 32000
$ws.Range("J44").Value = 32000
$ws.Range("L44").Value = 32000
$ws.Range("N44").Value = -33192

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H14").Value = 15000
$ws.Range("I14").Value = 15000
$ws.Range("K14").Value = 15000
$ws.Range("M14").Value = -14828
$ws.Range("H20").Value = 2985.5715
$ws.Range("I20").Value = 2799.8
$ws.Range("K20").Value = 2799.8
$ws.Range("M20").Value = -2573.8
$ws.Range("H136").Value = 7023.25
$ws.Range("I136").Value = 6466
$ws.Range("K136").Value = 19398
$ws.Range("M136").Value = -16848

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H21").Value = 34332.668
$ws.Range("I21").Value = 34332.668
$ws.Range("K21").Value = 34332.668
$ws.Range("M21").Value = -34097.668
$ws.Range("H30").Value = 25000
$ws.Range("I30").Value = 25000
$ws.Range("K30").Value = 25000
$ws.Range("M30").Value = -24893
$ws.Range("H35").Value = 34332.668
$ws.Range("I35").Value = 34332.668
$ws.Range("K35").Value = 34332.668
$ws.Range("M35").Value = -34042.668
$ws.Range("H126").Value = 5982.0586
$ws.Range("I126").Value = 997
$ws.Range("K126").Value = 2991
$ws.Range("M126").Value = -521
$ws.Range("H133").Value = 119298.8
$ws.Range("J133").Value = 119298.8
$ws.Range("L133").Value = 119298.8
$ws.Range("N133").Value = -129418.8
